$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. "Status" text changed from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it appears
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Re-create the hyperlinks on zh-cn / de-de sheets so that a new
#    hyperlink (Latest Target File) appears in column I, next to the
#    existing Source File Name hyperlink in column A.
# ---------------------------------------------------------------------------
$srcDisplay = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef.md"
$srcTarget  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fee78cd77b339453fe303c0ff3343e06aba388b9/e2e/32fa0c8e-2890-4fb7-90f9-7212dd3922ef.md"
$dupDisplay = "fffffc8c9a4e-3817-4304-8563-e63318d4b77d.md"
$dupTarget  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fee78cd77b339453fe303c0ff3343e06aba388b9/e2e/fffffc8c9a4e-3817-4304-8563-e63318d4b77d.md"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $srcTarget, "", "", $srcDisplay)
$ws2.Hyperlinks.Add($ws2.Range("I2"), $srcTarget, "", "", $srcDisplay)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $dupTarget, "", "", $dupDisplay)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $srcTarget, "", "", $srcDisplay)

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $srcTarget, "", "", $srcDisplay)
$ws3.Hyperlinks.Add($ws3.Range("I2"), $srcTarget, "", "", $srcDisplay)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $dupTarget, "", "", $dupDisplay)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $srcTarget, "", "", $srcDisplay)

# ---------------------------------------------------------------------------
# 3. Fill in "Latest Handback File" (column J) with the generated xliff
#    name for each locale.
# ---------------------------------------------------------------------------
$ws2.Range("J2").Value = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef.1b21172b5759f6b658c312369ecda99410f0a9d4.zh-cn.xlf"
$ws2.Range("J3").Value = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef.1b21172b5759f6b658c312369ecda99410f0a9d4.zh-cn.xlf"

$ws3.Range("J2").Value = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef.1b21172b5759f6b658c312369ecda99410f0a9d4.de-de.xlf"
$ws3.Range("J3").Value = "32fa0c8e-2890-4fb7-90f9-7212dd3922ef.1b21172b5759f6b658c312369ecda99410f0a9d4.de-de.xlf"

# ---------------------------------------------------------------------------
# 4. "Latest Handback DateTime" (column K):
#    zh-cn now has an actual handback time (was the zero/unset date).
#    de-de gets its own, later, handback time.
# ---------------------------------------------------------------------------
$ws2.Range("K2").Value = "2016-09-05 05:09:29"
$ws2.Range("K3").Value = "2016-09-05 05:09:29"

$ws3.Range("K2").Value = "2016-09-05 05:09:37"
$ws3.Range("K3").Value = "2016-09-05 05:09:37"

# ---------------------------------------------------------------------------
# 5. Widen the columns that now hold the longer strings.
#    (ColumnWidth is expressed in "characters"; Excel stores the column
#    definition in the file using a slightly different, pixel-rounded
#    unit, so the values below are chosen to serialize to the intended
#    stored widths of 30 and 40.)
# ---------------------------------------------------------------------------
$ws1.Range("E1:F1").ColumnWidth = 29.17

$ws2.Range("C1").ColumnWidth = 29.17
$ws2.Range("I1").ColumnWidth = 39.17
$ws2.Range("J1").ColumnWidth = 39.17

$ws3.Range("C1").ColumnWidth = 29.17
$ws3.Range("I1").ColumnWidth = 39.17
$ws3.Range("J1").ColumnWidth = 39.17
